# Actualización 11 de febrero de 2024 - Lap HP
# Se actualiza el repositorio con varios materiales.

$wb = $excel.ActiveWorkbook
$wsFaltas = $wb.Worksheets.Item("Faltas")
$wsConcentrado = $wb.Worksheets.Item("Concentrado")

# --- Faltas: add a new weekly attendance column (J) with its header date ---
$wsFaltas.Range("J1").Value = 45321
$wsFaltas.Range("J1").NumberFormat = "d-mmm"
$wsFaltas.Range("J1").HorizontalAlignment = -4108

# --- Faltas: fill in the missing absence marks for the existing weeks (H, I) ---
$wsFaltas.Range("H3").Value = 1
$wsFaltas.Range("H6").Value = 1
$wsFaltas.Range("H7").Value = 1
$wsFaltas.Range("I10").Value = 1
$wsFaltas.Range("I11").Value = 1
$wsFaltas.Range("H13").Value = 1

# --- Faltas: tighten up the weekly mark columns now that a new week was added ---
$wsFaltas.Range("F1:I13").ColumnWidth = 6.3
$wsFaltas.Range("J1:J13").ColumnWidth = 6.0

# --- Concentrado keeps its own selection, it just stops being the active tab ---
$wsConcentrado.Activate() | Out-Null
$wsConcentrado.Range("H8").Select() | Out-Null

# --- Selection / active sheet: Faltas becomes the active tab with K12 selected ---
$wsFaltas.Activate() | Out-Null
$wsFaltas.Range("K12").Select() | Out-Null
